# Assign Shift.xlsx - add a "Sheet2" lookup sheet holding the list of shift
# names, then wire it up as a dropdown (list data validation) on column C of
# Sheet1, which is where the ShiftName value is entered for each row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so tab order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$shifts = @(
    "Early Shift(06:00-14:30)",
    "Saturday Shift(06:00-11:00)",
    "Morning Shift(07:00-15:30)",
    "General Shift(08:00-16:30)",
    "Sat Shift(08:00-13:00)",
    "Day Shift(09:00-17:30)",
    "Extended Day Shift(10:00-18:30)",
    "Mid Shift(12:00-20:30)",
    "Afternoon Shift(14:00-22:30)",
    "Evening Shift(15:00-00:00)",
    "Late Evening Shift(16:00-01:00)",
    "Late Shift(17:00-02:00)",
    "Night Shift(18:00-03:00)",
    "Extended Night Shift(18:30-03:30)",
    "Late Night Shift(19:30-04:30)",
    "Overnight Shift(21:00-06:00)",
    "Midnight Shift(22:00-07:00)",
    "Weekly Off(00:00-00:00)"
)

for ($i = 0; $i -lt $shifts.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $shifts[$i]
}
$ws2.Columns.Item(1).ColumnWidth = 29.25  # renders as the ~30.14-char column width used in the template

# Make Sheet1 the active sheet again and select column C so the ShiftName
# column is the one the dropdown applies to.
$ws1.Activate()

$range = $ws1.Range("C1:C1048576")
$range.Validation.Add(3, 1, 1, "=Sheet2!`$A`$1:`$A`$18")
$range.Validation.IgnoreBlank = $true
$range.Validation.InCellDropdown = $true
$range.Validation.ShowInput = $true
$range.Validation.ShowError = $true

$null = $ws1.Range("C1").Select()
